# Insert a new data row at row 178 (shifts existing rows 178..279 down to 179..280),
# matching the target diff where a new "Brócoli" price record is added for
# Feria Lagunitas de Puerto Montt and all subsequent rows move down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 178; this pushes rows 178-279 to 179-280.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with the new record's data.
$ws.Range("A178").Value = 4
$ws.Range("B178").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value = "Los Lagos"
$ws.Range("D178").Value = 44596
$ws.Range("E178").Value = 10
$ws.Range("F178").Value = 100112023
$ws.Range("G178").Value = "Brócoli"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 600
$ws.Range("K178").Value = 1500
$ws.Range("L178").Value = 1500
$ws.Range("M178").Value = 1500
$ws.Range("N178").Value = "`$/unidad"
$ws.Range("O178").Value = "Región Metropolitana"
$ws.Range("P178").Value = 1500
$ws.Range("Q178").Value = 1
$ws.Range("R178").Value = "Hortaliza"
